$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The underlying Power Query refreshed and filled in the "Progress" value
# for the trial row in A14 ("REJOICE (MK-5909-003)"), which previously had
# no Progress figure, the same way it already did for other trial rows.
$ws.Range("B14").Value = 0
